$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The used range is A:AY. We copy row content in chunks that deliberately
# skip columns I, Y, AA and AT. Those columns are identical (either the
# empty placeholder cell, or the literal date string "2023-09-15") across
# every row touched below, so skipping them is a no-op for the final data
# -- but it also sidesteps two COM quirks that would otherwise create a
# spurious diff:
#   * assigning a date-like string ("2023-09-15") through Range.Value
#     triggers Excel's automatic text->date coercion, turning the cell
#     into a date serial number with a new style;
#   * assigning an empty string back through Range.Value drops the cell
#     entirely instead of keeping the original empty placeholder cell.

function Get-RowChunks($rowNum) {
    $a = $ws.Range("A" + $rowNum + ":H" + $rowNum).Value()
    $b = $ws.Range("J" + $rowNum + ":X" + $rowNum).Value()
    $c = $ws.Range("Z" + $rowNum + ":Z" + $rowNum).Value()
    $d = $ws.Range("AB" + $rowNum + ":AS" + $rowNum).Value()
    $e = $ws.Range("AU" + $rowNum + ":AY" + $rowNum).Value()
    $result = @($a, $b, $c, $d, $e)
    return $result
}

function Set-RowChunks($rowNum, $chunks) {
    $ws.Range("A" + $rowNum + ":H" + $rowNum).Value = $chunks[0]
    $ws.Range("J" + $rowNum + ":X" + $rowNum).Value = $chunks[1]
    $ws.Range("Z" + $rowNum + ":Z" + $rowNum).Value = $chunks[2]
    $ws.Range("AB" + $rowNum + ":AS" + $rowNum).Value = $chunks[3]
    $ws.Range("AU" + $rowNum + ":AY" + $rowNum).Value = $chunks[4]
}

# --- Rows 7-10: cyclic rotation ---
# New row7 = old row10, new row8 = old row7, new row9 = old row8, new row10 = old row9
$row7 = Get-RowChunks 7
$row8 = Get-RowChunks 8
$row9 = Get-RowChunks 9
$row10 = Get-RowChunks 10

Set-RowChunks 7 $row10
Set-RowChunks 8 $row7
Set-RowChunks 9 $row8
Set-RowChunks 10 $row9

# --- Rows 31 and 33: swap ---
$row31 = Get-RowChunks 31
$row33 = Get-RowChunks 33

Set-RowChunks 31 $row33
Set-RowChunks 33 $row31

# --- Rows 38 and 39: swap ---
$row38 = Get-RowChunks 38
$row39 = Get-RowChunks 39

Set-RowChunks 38 $row39
Set-RowChunks 39 $row38
